# Regenerate the "K" column (column G) of the save_data sheet.
# The data in this sheet is raw literal game-log data (no formulas), so the
# refreshed "K" (strikeouts) values are written directly as computed by the
# upstream data-regeneration script described in the commit message
# ("regen save_data to use K instead of Strike#, regen std/mean, calc and
# write s_vals").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$kValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 2
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    12 = 1
    13 = 0
    14 = 0
    15 = 1
    16 = 1
    18 = 0
    19 = 1
    20 = 2
    21 = 0
    22 = 0
    23 = 2
    24 = 1
    25 = 1
    26 = 1
    27 = 2
    28 = 2
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 1
    34 = 0
    35 = 2
    36 = 0
    37 = 0
    38 = 1
    39 = 0
    40 = 1
    41 = 2
    42 = 0
    43 = 4
    44 = 4
    45 = 0
    46 = 0
    47 = 1
    48 = 0
    49 = 0
    50 = 3
    51 = 1
    52 = 3
    53 = 1
    56 = 1
    57 = 1
    58 = 1
    59 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
